# repull data, push all data, mean calculation
# Update the dSF column (column F) with newly recalculated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 8
    3  = 5
    5  = -2
    7  = -3
    9  = -3
    10 = -2
    11 = -1
    12 = -2
    13 = -1
    14 = -6
    15 = 2
    19 = 1
    20 = -2
    21 = -1
    22 = 3
    23 = 3
    24 = -3
    25 = -2
    26 = -1
    27 = -1
    28 = 2
    29 = 2
    30 = -1
    31 = 2
    33 = 3
    35 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
